# New submission synced into the "JSS 3D" sheet: append one row (row 12)
# with the Google-Forms-style columns Timestamp / Full Name / Admission No / AI Score.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Range("A12").Value = "2026-02-12 07:39:40"
$ws.Range("B12").Value = "Abdullahi Tijjani Buji"
# Leading apostrophe forces this numeric-looking admission number to stay text,
# matching the other rows in column C (e.g. "38", "24", "7").
$ws.Range("C12").Value = "'19"
$ws.Range("D12").Value = 10
